$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 256-257 (existing rows 256-295 shift down to 258-297)
$ws.Rows("256:257").Insert()

# Row 256: new "Primera" quality record for 2021-10-05 (serial 44474)
$ws.Range("A256").Value = 3
$ws.Range("B256").Value = "Femacal de La Calera"
$ws.Range("C256").Value = "Coquimbo"
$ws.Range("D256").Value = 44474
$ws.Range("E256").Value = 5
$ws.Range("F256").Value = 100112037
$ws.Range("G256").Value = "Cebollín"
$ws.Range("H256").Value = "Sin especificar"
$ws.Range("I256").Value = "Primera"
$ws.Range("J256").Value = 280
$ws.Range("K256").Value = 2800
$ws.Range("L256").Value = 3000
$ws.Range("M256").Value = 2886
$ws.Range("N256").Value = "$/paquete 36 unidades"
$ws.Range("O256").Value = "Provincia de Quillota"
$ws.Range("P256").Value = 80
$ws.Range("Q256").Value = 36
$ws.Range("R256").Value = "Hortaliza"

# Row 257: new "Segunda" quality record for the same date
$ws.Range("A257").Value = 3
$ws.Range("B257").Value = "Femacal de La Calera"
$ws.Range("C257").Value = "Coquimbo"
$ws.Range("D257").Value = 44474
$ws.Range("E257").Value = 5
$ws.Range("F257").Value = 100112037
$ws.Range("G257").Value = "Cebollín"
$ws.Range("H257").Value = "Sin especificar"
$ws.Range("I257").Value = "Segunda"
$ws.Range("J257").Value = 120
$ws.Range("K257").Value = 2000
$ws.Range("L257").Value = 2000
$ws.Range("M257").Value = 2000
$ws.Range("N257").Value = "$/paquete 36 unidades"
$ws.Range("O257").Value = "Provincia de Quillota"
$ws.Range("P257").Value = 56
$ws.Range("Q257").Value = 36
$ws.Range("R257").Value = "Hortaliza"
